# Apply "repull data" updates to column F (dSF) for specific rows.
# These reflect re-pulled source values that now differ from column E (dS0).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    2  = -1
    3  = 1
    5  = 2
    10 = -2
    11 = -8
    13 = 1
    16 = -1
    18 = -1
    28 = -2
    31 = -11
    32 = 3
    39 = -7
    40 = -2
    41 = -2
    42 = -3
    48 = -1
    63 = -1
    67 = -4
    70 = -1
    71 = -1
    74 = -1
}

foreach ($row in $updates.Keys) {
    $ws.Range("F$row").Value = $updates[$row]
}
